$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$link1 = "https://www.genomeweb.com/sequencing/thermo-fisher-ngs-assay-gets-fda-ok-cdx-dizals-zegfrovy-and-solid-tumor-profiling"
$link2 = "https://www.360dx.com/sequencing/thermo-fisher-ngs-assay-gets-fda-ok-cdx-dizals-zegfrovy-and-solid-tumor-profiling"
$keyword = "CDx"
$title = "Thermo Fisher NGS Assay Gets FDA OK as CDx for Dizal's Zegfrovy and for Solid Tumor Profiling"

# Row 30
$ws.Range("A30").Value = $link1
$ws.Hyperlinks.Add($ws.Range("A30"), $link1)
$ws.Range("A30").Style = $ws.Range("A2").Style
$ws.Range("B30").Value = $keyword
$ws.Range("C30").Value = $title

# Row 31
$ws.Range("A31").Value = $link2
$ws.Hyperlinks.Add($ws.Range("A31"), $link2)
$ws.Range("A31").Style = $ws.Range("A2").Style
$ws.Range("B31").Value = $keyword
$ws.Range("C31").Value = $title
